$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: UNM / Unum Group / 14.47 / 2.943B / 3.066
$ws.Cells.Item(1, 1).Value = "UNM"
$ws.Cells.Item(1, 2).Value = "Unum Group"
$ws.Cells.Item(1, 3).Value = 14.47
$ws.Cells.Item(1, 4).Value = "2.943B"
$ws.Cells.Item(1, 5).Value = 3.0659999999999998

# Row 2: MET / MetLife Inc / 32.39 / 29.40B / 3.446
$ws.Cells.Item(2, 1).Value = "MET"
$ws.Cells.Item(2, 2).Value = "MetLife Inc"
$ws.Cells.Item(2, 3).Value = 32.39
$ws.Cells.Item(2, 4).Value = "29.40B"
$ws.Cells.Item(2, 5).Value = 3.4460000000000002

# Row 3: SYF / Synchrony Financial / 17.63 / 10.29B / 3.971
$ws.Cells.Item(3, 1).Value = "SYF"
$ws.Cells.Item(3, 2).Value = "Synchrony Financial"
$ws.Cells.Item(3, 3).Value = 17.63
$ws.Cells.Item(3, 4).Value = "10.29B"
$ws.Cells.Item(3, 5).Value = 3.9710000000000001

# E1 previously used its own 3-decimal numeric format; it now shares the
# same "0." format as the other value cells (numFmtId 164).
$ws.Cells.Item(1, 5).NumberFormat = "0."

# Row heights
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 60

# Column widths (closest values achievable given engine's pixel rounding)
$ws.Columns.Item(2).ColumnWidth = 7.33
$ws.Columns.Item(3).ColumnWidth = 33.67
$ws.Columns.Item(4).ColumnWidth = 23.5
$ws.Columns.Item(5).ColumnWidth = 9.83
$ws.Columns.Item(6).ColumnWidth = 11

# Selection
$ws.Range("C12").Select()
